$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("339÷7=48, 3", $true, $false, $false, $false, $false, $true, 1, $false, "963÷4=240, 3", 2) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("178÷9=19, 7", $true, $false, $false, $false, $false, $true, 1, $false, "628÷5=125, 3", 2) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("806÷6=134, 2", $true, $false, $false, $false, $false, $true, 1, $false, "207÷4=51, 3", 2) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("669÷5=133, 4", $true, $false, $false, $false, $false, $true, 1, $false, "296÷2=148, 0", 2) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("884÷2=442, 0", $true, $false, $false, $false, $false, $true, 1, $false, "158÷7=22, 4", 2) | Out-Null

$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("540÷5=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "800÷9=88, 8", 2) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("251÷4=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "434÷2=217, 0", 2) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("777÷5=155, 2", $true, $false, $false, $false, $false, $true, 1, $false, "460÷2=230, 0", 2) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("675÷6=112, 3", $true, $false, $false, $false, $false, $true, 1, $false, "992÷3=330, 2", 2) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("950÷6=158, 2", $true, $false, $false, $false, $false, $true, 1, $false, "849÷2=424, 1", 2) | Out-Null

$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("634÷5=126, 4", $true, $false, $false, $false, $false, $true, 1, $false, "755÷3=251, 2", 2) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("293÷6=48, 5", $true, $false, $false, $false, $false, $true, 1, $false, "862÷8=107, 6", 2) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("505÷5=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "940÷4=235, 0", 2) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("345÷5=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "138÷6=23, 0", 2) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("721÷3=240, 1", $true, $false, $false, $false, $false, $true, 1, $false, "205÷7=29, 2", 2) | Out-Null

$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("749÷5=149, 4", $true, $false, $false, $false, $false, $true, 1, $false, "986÷7=140, 6", 2) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("921÷8=115, 1", $true, $false, $false, $false, $false, $true, 1, $false, "806÷6=134, 2", 2) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("886÷8=110, 6", $true, $false, $false, $false, $false, $true, 1, $false, "695÷6=115, 5", 2) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("496÷9=55, 1", $true, $false, $false, $false, $false, $true, 1, $false, "111÷4=27, 3", 2) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("196÷9=21, 7", $true, $false, $false, $false, $false, $true, 1, $false, "315÷3=105, 0", 2) | Out-Null

$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("305÷8=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "722÷5=144, 2", 2) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("252÷7=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "847÷3=282, 1", 2) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("106÷5=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "881÷7=125, 6", 2) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("523÷6=87, 1", $true, $false, $false, $false, $false, $true, 1, $false, "639÷2=319, 1", 2) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("199÷6=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "851÷3=283, 2", 2) | Out-Null

Write-Host "Done applying replacements"

